$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the full data block (rows 2-24, columns A-T) BEFORE any writes,
# since the update re-orders rows (some rows are both a source and a
# destination), so we must not overwrite source data before it is read.
$firstRow = 2
$lastRow = 24
$lastCol = 20  # column T

$snapshot = $ws.Range($ws.Cells.Item($firstRow, 1), $ws.Cells.Item($lastRow, $lastCol)).Value2

# Mapping: new row number -> old row number whose data it should receive.
$mapping = @{
    2  = 7
    3  = 6
    4  = 12
    5  = 13
    6  = 15
    7  = 22
    8  = 19
    9  = 10
    10 = 14
    11 = 23
    12 = 24
    13 = 20
    14 = 4
    15 = 3
    16 = 8
    17 = 17
    18 = 11
    19 = 18
    20 = 5
    21 = 16
    22 = 21
    23 = 9
    24 = 2
}

foreach ($newRow in $mapping.Keys) {
    $oldRow = $mapping[$newRow]
    $srcIdx = $oldRow - $firstRow + 1
    for ($col = 1; $col -le $lastCol; $col++) {
        $value = $snapshot[$srcIdx, $col]
        $ws.Cells.Item($newRow, $col).Value = $value
    }
}
